$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 375.53845
$ws.Range("I33").Value = 414.9565
$ws.Range("K33").Value = 414.9565
$ws.Range("M33").Value = -185.9565
$ws.Range("H40").Value = 3449.9333
$ws.Range("I40").Value = 2122.6365
$ws.Range("J40").Value = 7100
$ws.Range("K40").Value = 2122.6365
$ws.Range("L40").Value = 7100
$ws.Range("M40").Value = -1947.6365
$ws.Range("N40").Value = -7450
$ws.Range("H41").Value = 281
$ws.Range("I41").Value = 336
$ws.Range("K41").Value = 336
$ws.Range("M41").Value = 104
$ws.Range("H62").Value = 9231.8
$ws.Range("I62").Value = 9231.8
$ws.Range("K62").Value = 9231.8
$ws.Range("M62").Value = -8607.8
$ws.Range("H65").Value = 9231.8
$ws.Range("I65").Value = 9231.8
$ws.Range("K65").Value = 46159
$ws.Range("M65").Value = -43039
$ws.Range("H69").Value = 45388.832
$ws.Range("I69").Value = 2183.4443
$ws.Range("K69").Value = 6550.3329
$ws.Range("M69").Value = -5676.3329
$ws.Range("H72").Value = 45388.832
$ws.Range("I72").Value = 2183.4443
$ws.Range("K72").Value = 19650.9987
$ws.Range("M72").Value = -15282.9987
$ws.Range("H96").Value = 940.6667
$ws.Range("I96").Value = 419.66666
$ws.Range("K96").Value = 1258.99998
$ws.Range("M96").Value = 114.0000199999999
$ws.Range("H98").Value = 2036.3889
$ws.Range("I98").Value = 962.6429
$ws.Range("K98").Value = 962.6429
$ws.Range("M98").Value = 535.3571
$ws.Range("H107").Value = 931.95654
$ws.Range("I107").Value = 952.1
$ws.Range("K107").Value = 952.1
$ws.Range("M107").Value = 967.9
$ws.Range("H113").Value = 3474.25
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 3613.4285
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 3613.4285
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -10121.4285
$ws.Range("H122").Value = 2036.3889
$ws.Range("I122").Value = 962.6429
$ws.Range("K122").Value = 2887.9287
$ws.Range("M122").Value = -437.9287000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1745.6875
$ws.Range("J2").Value = 1749.75
$ws.Range("L2").Value = 1749.75
$ws.Range("N2").Value = -1975.75
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28:N28").ClearContents()
$ws.Range("H32").Value = 11305.464
$ws.Range("J32").Value = 10037
$ws.Range("L32").Value = 10037
$ws.Range("N32").Value = -10611
$ws.Range("H61").Value = 3184.889
$ws.Range("I61").Value = 2847.875
$ws.Range("K61").Value = 2847.875
$ws.Range("M61").Value = -2635.875
$ws.Range("H97").Value = 4168.7144
$ws.Range("I97").Value = 696.2
$ws.Range("J97").Value = 12850
$ws.Range("K97").Value = 696.2
$ws.Range("L97").Value = 12850
$ws.Range("M97").Value = -200.2
$ws.Range("N97").Value = -13842
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99:N99").ClearContents()
$ws.Range("H116").Value = 1745.6875
$ws.Range("J116").Value = 1749.75
$ws.Range("L116").Value = 1749.75
$ws.Range("N116").Value = -6337.75
$ws.Range("H136").Value = 3184.889
$ws.Range("I136").Value = 2847.875
$ws.Range("K136").Value = 8543.625
$ws.Range("M136").Value = -5993.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1745.6875
$ws.Range("J3").Value = 1749.75
$ws.Range("L3").Value = 1749.75
$ws.Range("N3").Value = -1977.75
$ws.Range("H99").Value = 17759.916
$ws.Range("I99").Value = 21151.9
$ws.Range("K99").Value = 21151.9
$ws.Range("M99").Value = -19653.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2060.353
$ws.Range("I16").Value = 1534.909
$ws.Range("J16").Value = 3023.6667
$ws.Range("K16").Value = 1534.909
$ws.Range("L16").Value = 3023.6667
$ws.Range("M16").Value = -1247.909
$ws.Range("N16").Value = -3597.6667
$ws.Range("H86").Value = 4246
$ws.Range("J86").Value = 3029.3333
$ws.Range("L86").Value = 3029.3333
$ws.Range("N86").Value = -5275.3333
$ws.Range("H89").Value = 4246
$ws.Range("J89").Value = 3029.3333
$ws.Range("L89").Value = 15146.6665
$ws.Range("N89").Value = -26378.6665
$ws.Range("H113").Value = 2060.353
$ws.Range("I113").Value = 1534.909
$ws.Range("J113").Value = 3023.6667
$ws.Range("K113").Value = 1534.909
$ws.Range("L113").Value = 3023.6667
$ws.Range("M113").Value = 635.0909999999999
$ws.Range("N113").Value = -7363.6667
$ws.Range("H134").Value = 2638.625
$ws.Range("I134").Value = 2469.8
$ws.Range("K134").Value = 7409.400000000001
$ws.Range("M134").Value = -4874.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10636.272
$ws.Range("J34").Value = 10636.272
$ws.Range("L34").Value = 31908.816
$ws.Range("N34").Value = -32076.816
$ws.Range("H55").Value = 10247.5
$ws.Range("J55").Value = 13139.286
$ws.Range("L55").Value = 39417.858
$ws.Range("N55").Value = -39771.858
$ws.Range("H121").Value = 2393.9
$ws.Range("J121").Value = 2683.75
$ws.Range("L121").Value = 8051.25
$ws.Range("N121").Value = -10671.25
$ws.Range("H132").Value = 2356.8333
$ws.Range("I132").Value = 1110.25
$ws.Range("J132").Value = 4850
$ws.Range("K132").Value = 9992.25
$ws.Range("L132").Value = 43650
$ws.Range("M132").Value = -7462.25
$ws.Range("N132").Value = -48710

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1906.0435
$ws.Range("I93").Value = 1785.2142
$ws.Range("J93").Value = 2094
$ws.Range("K93").Value = 1785.2142
$ws.Range("L93").Value = 2094
$ws.Range("M93").Value = -537.2141999999999
$ws.Range("N93").Value = -4590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 25951.334
$ws.Range("I41").Value = 28943
$ws.Range("J41").Value = 19968
$ws.Range("K41").Value = 28943
$ws.Range("L41").Value = 19968
$ws.Range("M41").Value = -28553
$ws.Range("N41").Value = -20748
$ws.Range("H81").Value = 1707.625
$ws.Range("I81").Value = 1634.1428
$ws.Range("J81").Value = 2222
$ws.Range("K81").Value = 3268.2856
$ws.Range("L81").Value = 4444
$ws.Range("M81").Value = -2207.2856
$ws.Range("N81").Value = -6566
$ws.Range("H84").Value = 1707.625
$ws.Range("I84").Value = 1634.1428
$ws.Range("J84").Value = 2222
$ws.Range("K84").Value = 16341.428
$ws.Range("L84").Value = 22220
$ws.Range("M84").Value = -11037.428
$ws.Range("N84").Value = -32828
$ws.Range("H113").Value = 474.5
$ws.Range("I113").Value = 421.2857
$ws.Range("J113").Value = 515.8889
$ws.Range("K113").Value = 1263.8571
$ws.Range("L113").Value = 1547.6667
$ws.Range("M113").Value = 906.1428999999998
$ws.Range("N113").Value = -5887.6667
$ws.Range("H132").Value = 2285.8
$ws.Range("I132").Value = 1223
$ws.Range("K132").Value = 3669
$ws.Range("M132").Value = -1139
$ws.Range("H136").Value = 5793.95
$ws.Range("I136").Value = 5051.706
$ws.Range("K136").Value = 15155.118
$ws.Range("M136").Value = -12605.118
